$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell, even if it looks like a
# number/currency (e.g. "$879.38"), so it is stored as text rather than
# being auto-converted to a formatted number.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- GLASS section relabeling (rows 13, 15, 16) ---
$ws.Range("F13").Value = "Glass Area"
$ws.Range("F15").Value = "83.75 sqft"
Set-TextValue $ws.Range("F16") "$879.38"

# --- New FABRICATION section ---
# Old GRAND TOTAL (rows 19-20) is replaced by the FABRICATION block, and
# GRAND TOTAL moves down to rows 25-26.
$ws.Range("E19").Value = "FABRICATION"
$ws.Range("F19").Value = "Joints Fabrication Labor"

$ws.Range("E20").Value = "Part Number"
$ws.Range("F20").Value = "N/A"

$ws.Range("E21").Value = "Quantity"
$ws.Range("F21").Value = "18 joints"

$ws.Range("E22").Value = "Price"
Set-TextValue $ws.Range("F22") "$270.00"

# Rows 23 and 24 remain blank spacer rows.

$ws.Range("E25").Value = "GRAND TOTAL"
Set-TextValue $ws.Range("E26") "$3303.86"
